$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (MuSCs and Resolving-Mac -> ECs/Resolving-Mac rows merged away)
$ws.Range("A6:T9").EntireRow.Delete() | Out-Null

# Row 2: ECs -> Resolving-Mac
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.013740666666667
$ws.Range("H2").Value = 6.041221999999999
$ws.Range("I2").Value = 0.3805515268368102
$ws.Range("J2").Value = 0.3805515268368102
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.569632
$ws.Range("N2").Value = 67.708896
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 45.44938579010133
$ws.Range("R2").Value = 409.044472110912
$ws.Range("S2").Value = 0.3805515268368102
$ws.Range("T2").Value = 0.3805515268368102

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8431363333333334
$ws.Range("H3").Value = 2.529409
$ws.Range("I3").Value = 0.1593337336295156
$ws.Range("J3").Value = 0.1593337336295156
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 22.569632
$ws.Range("N3").Value = 67.708896
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 19.02927676916267
$ws.Range("R3").Value = 171.263490922464
$ws.Range("S3").Value = 0.1593337336295156
$ws.Range("T3").Value = 0.1593337336295156

# Row 4: FAPs -> MuSCs, target ECs -> Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 1.000022
$ws.Range("H4").Value = 3.000066
$ws.Range("I4").Value = 0.1889815830160193
$ws.Range("J4").Value = 0.1889815830160193
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.569632
$ws.Range("N4").Value = 67.708896
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 22.570128531904
$ws.Range("R4").Value = 203.131156787136
$ws.Range("S4").Value = 0.1889815830160193
$ws.Range("T4").Value = 0.1889815830160193

# Row 5: FAPs -> Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 1.434738333333333
$ws.Range("H5").Value = 4.304214999999999
$ws.Range("I5").Value = 0.271133156517655
$ws.Range("J5").Value = 0.271133156517655
$ws.Range("M5").Value = 22.569632
$ws.Range("N5").Value = 67.708896
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 32.38151619962666
$ws.Range("R5").Value = 291.43364579664
$ws.Range("S5").Value = 0.271133156517655
$ws.Range("T5").Value = 0.271133156517655
